$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    if ($val -match '^-?[0-9]*\.?[0-9]+$') {
        $sheet.Range($addr).Value = "'" + $val
    } else {
        $sheet.Range($addr).Value = $val
    }
}

Set-TextValue $ws "D2" '28.422.56'
$ws.Range("E2").Value = '  +3.41%  '

Set-TextValue $ws "D3" '1.867.85'
$ws.Range("E3").Value = '  +2.00%  '

$ws.Range("E4").Value = '  -0.37%  '

Set-TextValue $ws "D5" '337.45'
$ws.Range("E5").Value = '  +1.88%  '

$ws.Range("E6").Value = '  -0.32%  '

Set-TextValue $ws "D7" '0.4711'
$ws.Range("E7").Value = '  +2.34%  '

$ws.Range("E8").Value = '  +3.78%  '

Set-TextValue $ws "D9" '47.55'
$ws.Range("E9").Value = '  +2.17%  '

Set-TextValue $ws "D10" '0.08012'
$ws.Range("E10").Value = '  +1.22%  '

Set-TextValue $ws "D11" '0.9995'
$ws.Range("E11").Value = '  +2.85%  '

Set-TextValue $ws "D12" '21.99'
$ws.Range("E12").Value = '  +4.22%  '

Set-TextValue $ws "D13" '6.037'
$ws.Range("E13").Value = '  +2.60%  '

Set-TextValue $ws "D14" '1.860.37'
$ws.Range("E14").Value = '  +1.93%  '

Set-TextValue $ws "D15" '7.257'
$ws.Range("E15").Value = '  +2.72%  '

Set-TextValue $ws "D16" '90.57'
$ws.Range("E16").Value = '  +2.71%  '

$ws.Range("E17").Value = '  -0.28%  '

Set-TextValue $ws "D18" '0.00001042'
$ws.Range("E18").Value = '  +1.01%  '

Set-TextValue $ws "D19" '0.06654'
$ws.Range("E19").Value = '  +0.06%  '

Set-TextValue $ws "D20" '17.57'
$ws.Range("E20").Value = '  +1.88%  '

$ws.Range("E21").Value = '  -0.25%  '

Set-TextValue $ws "D22" '28.448.56'
$ws.Range("E22").Value = '  +3.56%  '

Set-TextValue $ws "D23" '5.479'
$ws.Range("E23").Value = '  +2.58%  '

$ws.Range("E24").Value = '  +2.20%  '

Set-TextValue $ws "D25" '2.271'
$ws.Range("E25").Value = '  -1.63%  '

Set-TextValue $ws "D26" '2.085.02'
$ws.Range("E26").Value = '  +1.75%  '

Set-TextValue $ws "D27" '160.49'
$ws.Range("E27").Value = '  +1.97%  '

$ws.Range("E28").Value = '  +1.60%  '

Set-TextValue $ws "D29" '2.119'
$ws.Range("E29").Value = '  +2.68%  '

Set-TextValue $ws "D30" '5.484'
$ws.Range("E30").Value = '  +4.18%  '

Set-TextValue $ws "D31" '119.76'

Set-TextValue $ws "D32" '0.9662'
$ws.Range("E32").Value = '  +1.37%  '

Set-TextValue $ws "D33" '0.09517'
$ws.Range("E33").Value = '  +2.32%  '

Set-TextValue $ws "D34" '3.588'
$ws.Range("E34").Value = '  +0.55%  '

Set-TextValue $ws "D35" '1.380'
$ws.Range("E35").Value = '  +4.71%  '

$ws.Range("E36").Value = '  +2.07%  '

Set-TextValue $ws "D37" '0.06118'
$ws.Range("E37").Value = '  +2.99%  '

Set-TextValue $ws "D38" '0.02250'
$ws.Range("E38").Value = '  +2.21%  '

Set-TextValue $ws "D39" '8.315'
$ws.Range("E39").Value = '  +3.39%  '

$ws.Range("E40").Value = '  +2.24%  '

Set-TextValue $ws "D41" '0.5936'
$ws.Range("E41").Value = '  +2.44%  '

$ws.Range("E42").Value = '  -0.21%  '

$ws.Range("E43").Value = '  +1.84%  '

Set-TextValue $ws "D44" '10.33'
$ws.Range("E44").Value = '  +2.94%  '

$ws.Range("E45").Value = '  +0.28%  '

Set-TextValue $ws "D46" '0.5560'
$ws.Range("E46").Value = '  +1.19%  '

Set-TextValue $ws "D47" '12.14'
$ws.Range("E47").Value = '  +1.19%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws "D48" '1.955'
$ws.Range("E48").Value = '  +4.42%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws "D49" '0.07206'
$ws.Range("E49").Value = '  +8.36%  '

Set-TextValue $ws "D50" '2.067'
$ws.Range("E50").Value = '  +13.03%  '

Set-TextValue $ws "D51" '111.99'
$ws.Range("E51").Value = '  +1.58%  '
